# Insert two new weekly price records at the top of the price-history block
# (rows 177-178), shifting every existing record below down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 177:178 - everything currently at row 177 and
# below moves down by 2 (e.g. old row 177 -> new row 179, old row 274 ->
# new row 276), matching the target dimension A1:R276.
$ws.Rows("177:178").Insert()

# New row 177: Primera, Región de La Araucanía
$ws.Range("A177").Value = 11
$ws.Range("B177").Value = "Vega Monumental Concepción"
$ws.Range("C177").Value = "Bíobío"
$ws.Range("D177").Value = 44824
$ws.Range("E177").Value = 8
$ws.Range("F177").Value = 100114013
$ws.Range("G177").Value = "Zanahoria"
$ws.Range("H177").Value = "Sin especificar"
$ws.Range("I177").Value = "Primera"
$ws.Range("J177").Value = 800
$ws.Range("K177").Value = 10000
$ws.Range("L177").Value = 11000
$ws.Range("M177").Value = 10500
$ws.Range("N177").Value = "$/saco 20 kilos"
$ws.Range("O177").Value = "Región de La Araucanía"
$ws.Range("P177").Value = 525
$ws.Range("Q177").Value = 20
$ws.Range("R177").Value = "Hortaliza"

# New row 178: Segunda, Región de La Araucanía
$ws.Range("A178").Value = 11
$ws.Range("B178").Value = "Vega Monumental Concepción"
$ws.Range("C178").Value = "Bíobío"
$ws.Range("D178").Value = 44824
$ws.Range("E178").Value = 8
$ws.Range("F178").Value = 100114013
$ws.Range("G178").Value = "Zanahoria"
$ws.Range("H178").Value = "Sin especificar"
$ws.Range("I178").Value = "Segunda"
$ws.Range("J178").Value = 400
$ws.Range("K178").Value = 8000
$ws.Range("L178").Value = 8000
$ws.Range("M178").Value = 8000
$ws.Range("N178").Value = "$/saco 20 kilos"
$ws.Range("O178").Value = "Región de La Araucanía"
$ws.Range("P178").Value = 400
$ws.Range("Q178").Value = 20
$ws.Range("R178").Value = "Hortaliza"
